# Update "想去人数" (interest count) values on the 展览 and 全部类型 sheets
# to match the freshly generated gh-pages data output.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates
$wsExhibit.Range("F5").Value  = 915
$wsExhibit.Range("F7").Value  = 1882
$wsExhibit.Range("F10").Value = 60
$wsExhibit.Range("F12").Value = 511
$wsExhibit.Range("F13").Value = 188
$wsExhibit.Range("F18").Value = 8788
$wsExhibit.Range("F20").Value = 6851
$wsExhibit.Range("F21").Value = 11119
$wsExhibit.Range("F27").Value = 2446
$wsExhibit.Range("F30").Value = 2293
$wsExhibit.Range("F31").Value = 437
$wsExhibit.Range("F33").Value = 4463
$wsExhibit.Range("F34").Value = 644

# 全部类型 sheet updates (same events mirrored into the combined sheet)
$wsAll.Range("F9").Value  = 915
$wsAll.Range("F11").Value = 1882
$wsAll.Range("F17").Value = 511
$wsAll.Range("F18").Value = 188
$wsAll.Range("F23").Value = 8788
$wsAll.Range("F25").Value = 6851
$wsAll.Range("F26").Value = 11119
$wsAll.Range("F38").Value = 4463
